$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Trening" header in F1, copying header style from A1
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Replace column A values with Excel date serials, and assign datetime number format
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A2").Copy()
$ws.Range("A3:A13").PasteSpecial(-4122)

# Fill in the full dataset (12 rows), columns A-F
$ws.Cells.Item(2, 1).Value = 45684.59250381945
$ws.Cells.Item(2, 2).Value = 592.3
$ws.Cells.Item(2, 3).Value = 11.68
$ws.Cells.Item(2, 4).Value = 2.128043089594159
$ws.Cells.Item(2, 5).Value = "10-15"
$ws.Cells.Item(2, 6).Value = "Duża Gra"

$ws.Cells.Item(3, 1).Value = 45684.59320289352
$ws.Cells.Item(3, 2).Value = 652.7
$ws.Cells.Item(3, 3).Value = 11.85
$ws.Cells.Item(3, 4).Value = 2.27735161781311
$ws.Cells.Item(3, 5).Value = "10-15"
$ws.Cells.Item(3, 6).Value = "Duża Gra"

$ws.Cells.Item(4, 1).Value = 45684.59377233796
$ws.Cells.Item(4, 2).Value = 701.9
$ws.Cells.Item(4, 3).Value = 10.57
$ws.Cells.Item(4, 4).Value = 1.798816118921554
$ws.Cells.Item(4, 5).Value = "10-15"
$ws.Cells.Item(4, 6).Value = "Duża Gra"

$ws.Cells.Item(5, 1).Value = 45684.59250150463
$ws.Cells.Item(5, 2).Value = 592.1
$ws.Cells.Item(5, 3).Value = 9.9
$ws.Cells.Item(5, 4).Value = 2.062567251069203
$ws.Cells.Item(5, 5).Value = "5-10"
$ws.Cells.Item(5, 6).Value = "Duża Gra"

$ws.Cells.Item(6, 1).Value = 45684.5931994213
$ws.Cells.Item(6, 2).Value = 652.4
$ws.Cells.Item(6, 3).Value = 9.1
$ws.Cells.Item(6, 4).Value = 2.064195905412946
$ws.Cells.Item(6, 5).Value = "5-10"
$ws.Cells.Item(6, 6).Value = "Duża Gra"

$ws.Cells.Item(7, 1).Value = 45684.59377118055
$ws.Cells.Item(7, 2).Value = 701.8
$ws.Cells.Item(7, 3).Value = 9.88
$ws.Cells.Item(7, 4).Value = 1.756729228155956
$ws.Cells.Item(7, 5).Value = "5-10"
$ws.Cells.Item(7, 6).Value = "Duża Gra"

$ws.Cells.Item(8, 1).Value = 45684.59875844907
$ws.Cells.Item(8, 2).Value = 1132.7
$ws.Cells.Item(8, 3).Value = 14.37
$ws.Cells.Item(8, 4).Value = 3.623019116265433
$ws.Cells.Item(8, 5).Value = "10-15"
$ws.Cells.Item(8, 6).Value = "Mała Gra"

$ws.Cells.Item(9, 1).Value = 45684.60087418981
$ws.Cells.Item(9, 2).Value = 1315.5
$ws.Cells.Item(9, 3).Value = 11.93
$ws.Cells.Item(9, 4).Value = 3.698765754699708
$ws.Cells.Item(9, 5).Value = "10-15"
$ws.Cells.Item(9, 6).Value = "Mała Gra"

$ws.Cells.Item(10, 1).Value = 45684.60149340278
$ws.Cells.Item(10, 2).Value = 1369
$ws.Cells.Item(10, 3).Value = 14.81
$ws.Cells.Item(10, 4).Value = 3.94720697402954
$ws.Cells.Item(10, 5).Value = "10-15"
$ws.Cells.Item(10, 6).Value = "Mała Gra"

$ws.Cells.Item(11, 1).Value = 45684.59767280093
$ws.Cells.Item(11, 2).Value = 1038.9
$ws.Cells.Item(11, 3).Value = 9.76
$ws.Cells.Item(11, 4).Value = 2.859877824783324
$ws.Cells.Item(11, 5).Value = "5-10"
$ws.Cells.Item(11, 6).Value = "Mała Gra"

$ws.Cells.Item(12, 1).Value = 45684.59814386574
$ws.Cells.Item(12, 2).Value = 1079.6
$ws.Cells.Item(12, 3).Value = 9.43
$ws.Cells.Item(12, 4).Value = 3.194271530423845
$ws.Cells.Item(12, 5).Value = "5-10"
$ws.Cells.Item(12, 6).Value = "Mała Gra"

$ws.Cells.Item(13, 1).Value = 45684.600871875
$ws.Cells.Item(13, 2).Value = 1315.3
$ws.Cells.Item(13, 3).Value = 8.79
$ws.Cells.Item(13, 4).Value = 3.504503079823086
$ws.Cells.Item(13, 5).Value = "5-10"
$ws.Cells.Item(13, 6).Value = "Mała Gra"
